# Iteration4: Code on CodeSandbox. Paperwork done
# Adds 5 new defect-log rows (16-20) to Sheet1, reusing the formatting of
# the existing row 21, and updates the active-sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- Row 22 (Defect #16) ------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial($xlPasteFormats)
$ws.Range("A22").Value = 43583
$ws.Range("B22").Value = 16
$ws.Range("C22").Value = 80
$ws.Range("D22").Value = "Code"
$ws.Range("E22").Value = "Compile"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = '[Vue warn]: Property or method "tryLower" is not defined on the instance but referenced during render. Make sure that this property is reactive, either in the data option, or for class-based components, by initializing the property. See: https://vuejs.org/v2/guide/reactivity.html#Declaring-Reactive-Properties. - Failed to declare function in props'
$ws.Rows.Item(22).RowHeight = 75

# --- Row 23 (Defect #17) ------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A23").PasteSpecial($xlPasteFormats)
$ws.Range("A23").Value = 43583
$ws.Range("B23").Value = 17
$ws.Range("C23").Value = 20
$ws.Range("D23").Value = "Code"
$ws.Range("E23").Value = "Compile"
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = "SyntaxError: /src/App.vue: Unterminated string constant (223:11)"

# --- Row 24 (Defect #18) ------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A24").PasteSpecial($xlPasteFormats)
$ws.Range("A24").Value = 43583
$ws.Range("B24").Value = 18
$ws.Range("C24").Value = 70
$ws.Range("D24").Value = "Code"
$ws.Range("E24").Value = "Compile"
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 18
$ws.Range("H24").Value = "Random number would remain the same throughout testing, as well as highestNumber/lowestNumber. Added reset function to clear everything"
$ws.Rows.Item(24).RowHeight = 30

# --- Row 25 (Defect #19) ------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A25").PasteSpecial($xlPasteFormats)
$ws.Range("A25").Value = 43583
$ws.Range("B25").Value = 19
$ws.Range("C25").Value = 70
$ws.Range("D25").Value = "Code"
$ws.Range("E25").Value = "Compile"
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = '[Vue warn]: Invalid prop: type check failed for prop "currentStatement". Expected String with value "78", got Number with value 78. - .toString() added'
$ws.Rows.Item(25).RowHeight = 30

# --- Row 26 (Defect #20) ------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A26").PasteSpecial($xlPasteFormats)
$ws.Range("A26").Value = 43583
$ws.Range("B26").Value = 20
$ws.Range("C26").Value = 70
$ws.Range("D26").Value = "Code"
$ws.Range("E26").Value = "Compile"
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = '[Vue warn]: Invalid prop: type check failed for prop "currentStatement". Expected String with value "82", got Number with value 82. - .toString() added'
$ws.Rows.Item(26).RowHeight = 30

$ws.Application.CutCopyMode = $false

# --- Sheet view / selection ---------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I22").Select()
